# Auto-generated Excel COM-interop script
# Applies updated market-price values to the Excalibur_Profits workbook sheets
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 970.19147
$ws.Range("J17").Value = 990.2273
$ws.Range("L17").Value = 2970.6819
$ws.Range("N17").Value = -3306.6819
$ws.Range("H74").Value = 7088.7036
$ws.Range("I74").Value = 3798.7
$ws.Range("K74").Value = 3798.7
$ws.Range("M74").Value = -2862.7
$ws.Range("H76").Value = 3763.3635
$ws.Range("I76").Value = 3762.25
$ws.Range("K76").Value = 3762.25
$ws.Range("M76").Value = -3447.25
$ws.Range("H77").Value = 7088.7036
$ws.Range("I77").Value = 3798.7
$ws.Range("K77").Value = 18993.5
$ws.Range("M77").Value = -14313.5
$ws.Range("H79").Value = 3763.3635
$ws.Range("I79").Value = 3762.25
$ws.Range("K79").Value = 3762.25
$ws.Range("M79").Value = -2670.25
$ws.Range("H86").Value = 2600.5
$ws.Range("I86").Value = 2128.1428
$ws.Range("K86").Value = 2128.1428
$ws.Range("M86").Value = -1005.1428
$ws.Range("H89").Value = 2600.5
$ws.Range("I89").Value = 2128.1428
$ws.Range("K89").Value = 10640.714
$ws.Range("M89").Value = -5024.714
$ws.Range("H98").Value = 1072.7778
$ws.Range("I98").Value = 1068
$ws.Range("J98").Value = 1111
$ws.Range("K98").Value = 1068
$ws.Range("L98").Value = 1111
$ws.Range("M98").Value = 430
$ws.Range("N98").Value = -4107
$ws.Range("H122").Value = 1072.7778
$ws.Range("I122").Value = 1068
$ws.Range("J122").Value = 1111
$ws.Range("K122").Value = 3204
$ws.Range("L122").Value = 3333
$ws.Range("M122").Value = -754
$ws.Range("N122").Value = -8233
$ws.Range("H132").Value = 105305.914
$ws.Range("I132").Value = 132136.89
$ws.Range("K132").Value = 396410.67
$ws.Range("M132").Value = -393880.67
$ws.Range("H137").Value = 1160443.8
$ws.Range("I137").Value = 1074.4375
$ws.Range("K137").Value = 3223.3125
$ws.Range("M137").Value = -673.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2015.6364
$ws.Range("I33").Value = 2015.6364
$ws.Range("K33").Value = 2015.6364
$ws.Range("M33").Value = -1686.6364
$ws.Range("H61").Value = 1589188.9
$ws.Range("I61").Value = 1755998.2
$ws.Range("K61").Value = 1755998.2
$ws.Range("M61").Value = -1755786.2
$ws.Range("H74").Value = 2805.7
$ws.Range("I74").Value = 916.9091
$ws.Range("K74").Value = 916.9091
$ws.Range("M74").Value = -42.90909999999997
$ws.Range("H77").Value = 2805.7
$ws.Range("I77").Value = 916.9091
$ws.Range("K77").Value = 4584.5455
$ws.Range("M77").Value = -216.5455000000002
$ws.Range("H86").Value = 80235.5
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18814
$ws.Range("H88").Value = 1802.2222
$ws.Range("I88").Value = 3101
$ws.Range("J88").Value = 1431.1428
$ws.Range("K88").Value = 3101
$ws.Range("L88").Value = 1431.1428
$ws.Range("M88").Value = -2695
$ws.Range("N88").Value = -2243.1428
$ws.Range("H89").Value = 80235.5
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 60000
$ws.Range("M89").Value = -54072
$ws.Range("H91").Value = 1802.2222
$ws.Range("I91").Value = 3101
$ws.Range("J91").Value = 1431.1428
$ws.Range("K91").Value = 3101
$ws.Range("L91").Value = 1431.1428
$ws.Range("M91").Value = -1697
$ws.Range("N91").Value = -4239.1428
$ws.Range("H132").Value = 915202.9
$ws.Range("I132").Value = 1583155.9
$ws.Range("J132").Value = 4357.727
$ws.Range("K132").Value = 4749467.699999999
$ws.Range("L132").Value = 13073.181
$ws.Range("M132").Value = -4746937.699999999
$ws.Range("N132").Value = -18133.181
$ws.Range("H136").Value = 1589188.9
$ws.Range("I136").Value = 1755998.2
$ws.Range("K136").Value = 5267994.6
$ws.Range("M136").Value = -5265444.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1032724.3
$ws.Range("I58").Value = 2470788.8
$ws.Range("J58").Value = 5535.4287
$ws.Range("K58").Value = 2470788.8
$ws.Range("L58").Value = 5535.4287
$ws.Range("M58").Value = -2470585.8
$ws.Range("N58").Value = -5941.4287
$ws.Range("H99").Value = 1990.091
$ws.Range("I99").Value = 1802.5
$ws.Range("K99").Value = 1802.5
$ws.Range("M99").Value = -304.5
$ws.Range("H126").Value = 1990.091
$ws.Range("I126").Value = 1802.5
$ws.Range("K126").Value = 5407.5
$ws.Range("M126").Value = -2937.5
$ws.Range("H132").Value = 23378050
$ws.Range("I132").Value = 27029822
$ws.Range("J132").Value = 858782.7
$ws.Range("K132").Value = 81089466
$ws.Range("L132").Value = 2576348.1
$ws.Range("M132").Value = -81086936
$ws.Range("N132").Value = -2581408.1
$ws.Range("H134").Value = 3268536
$ws.Range("I134").Value = 9401.764999999999
$ws.Range("J134").Value = 12502750
$ws.Range("K134").Value = 28205.295
$ws.Range("L134").Value = 37508250
$ws.Range("M134").Value = -25670.295
$ws.Range("N134").Value = -37513320
$ws.Range("H136").Value = 1032724.3
$ws.Range("I136").Value = 2470788.8
$ws.Range("J136").Value = 5535.4287
$ws.Range("K136").Value = 7412366.399999999
$ws.Range("L136").Value = 16606.2861
$ws.Range("M136").Value = -7409816.399999999
$ws.Range("N136").Value = -21706.2861

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7362
$ws.Range("I75").Value = 1463
$ws.Range("J75").Value = 9983.777
$ws.Range("K75").Value = 4389
$ws.Range("L75").Value = 29951.331
$ws.Range("M75").Value = -3391
$ws.Range("N75").Value = -31947.331
$ws.Range("H78").Value = 7362
$ws.Range("I78").Value = 1463
$ws.Range("J78").Value = 9983.777
$ws.Range("K78").Value = 13167
$ws.Range("L78").Value = 89853.993
$ws.Range("M78").Value = -8175
$ws.Range("N78").Value = -99837.993
$ws.Range("H138").Value = 6004.4287
$ws.Range("I138").Value = 2998
$ws.Range("J138").Value = 6505.5
$ws.Range("K138").Value = 8994
$ws.Range("L138").Value = 19516.5
$ws.Range("M138").Value = -3854
$ws.Range("N138").Value = -29796.5
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 368245
$ws.Range("I80").Value = 396279.22
$ws.Range("J80").Value = 3800
$ws.Range("K80").Value = 396279.22
$ws.Range("L80").Value = 3800
$ws.Range("M80").Value = -395281.22
$ws.Range("N80").Value = -5796
$ws.Range("H83").Value = 368245
$ws.Range("I83").Value = 396279.22
$ws.Range("J83").Value = 3800
$ws.Range("K83").Value = 1981396.1
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = -1976404.1
$ws.Range("N83").Value = -28984
$ws.Range("H97").Value = 5331.231
$ws.Range("I97").Value = 5540.6
$ws.Range("K97").Value = 5540.6
$ws.Range("M97").Value = -5044.6
$ws.Range("H102").Value = 3262.652
$ws.Range("I102").Value = 2745.9512
$ws.Range("K102").Value = 2745.9512
$ws.Range("M102").Value = -1123.9512
$ws.Range("H107").Value = 36800
$ws.Range("I107").Value = 59466.168
$ws.Range("J107").Value = 2800.75
$ws.Range("K107").Value = 59466.168
$ws.Range("L107").Value = 2800.75
$ws.Range("M107").Value = -57546.168
$ws.Range("N107").Value = -6640.75
$ws.Range("H126").Value = 928815.2
$ws.Range("I126").Value = 1517441.6
$ws.Range("J126").Value = 3830.7144
$ws.Range("K126").Value = 4552324.800000001
$ws.Range("L126").Value = 11492.1432
$ws.Range("M126").Value = -4549854.800000001
$ws.Range("N126").Value = -16432.1432
$ws.Range("H132").Value = 36153056
$ws.Range("I132").Value = 44010540
$ws.Range("J132").Value = 8638.6
$ws.Range("K132").Value = 132031620
$ws.Range("L132").Value = 25915.8
$ws.Range("M132").Value = -132029090
$ws.Range("N132").Value = -30975.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3787.24
$ws.Range("I7").Value = 3613.3809
$ws.Range("K7").Value = 3613.3809
$ws.Range("M7").Value = -3501.3809
$ws.Range("H82").Value = 1270.6364
$ws.Range("J82").Value = 1307.8
$ws.Range("L82").Value = 1307.8
$ws.Range("N82").Value = -2029.8
$ws.Range("H85").Value = 1270.6364
$ws.Range("J85").Value = 1307.8
$ws.Range("L85").Value = 1307.8
$ws.Range("N85").Value = -3803.8
$ws.Range("H126").Value = 3787.24
$ws.Range("I126").Value = 3613.3809
$ws.Range("K126").Value = 10840.1427
$ws.Range("M126").Value = -8370.1427

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1999.3334
$ws.Range("I107").Value = 918.8461
$ws.Range("J107").Value = 3755.125
$ws.Range("K107").Value = 2756.5383
$ws.Range("L107").Value = 11265.375
$ws.Range("M107").Value = -836.5383000000002
$ws.Range("N107").Value = -15105.375
$ws.Range("H132").Value = 7456668
$ws.Range("I132").Value = 8386085
$ws.Range("K132").Value = 25158255
$ws.Range("M132").Value = -25155725
